# Updated symbol list on Sun Dec 11 21:49:33 UTC 2022 with GitHub Actions
#
# All data cells in this sheet are stored as text (inline strings), even
# the numeric-looking "Price" column. Writing a numeric-looking string via
# .Value would normally be auto-coerced to a real number by Excel, so the
# price-column updates are written with the cell temporarily forced to
# Text format ("@") and then the formatting is cleared again so no stray
# style survives the round-trip - this keeps the stored cell type as Text
# (matching the original inlineStr cells) while leaving cell formatting
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# --- Price ("D") column updates -------------------------------------------
Set-TextValue "D2"  "288.86"
Set-TextValue "D3"  "21.23"
Set-TextValue "D4"  "6.452"
Set-TextValue "D5"  "0.06385"
Set-TextValue "D6"  "3.604"
Set-TextValue "D7"  "1.578"
Set-TextValue "D8"  "6.591"
Set-TextValue "D9"  "0.8260"
Set-TextValue "D10" "0.01426"
Set-TextValue "D11" "0.1687"
Set-TextValue "D12" "0.08804"
Set-TextValue "D13" "0.03682"
Set-TextValue "D14" "0.03206"
Set-TextValue "D15" "0.09195"
Set-TextValue "D16" "3.715"
Set-TextValue "D17" "0.001643"
Set-TextValue "D18" "0.04742"
Set-TextValue "D19" "0.006128"
Set-TextValue "D20" "0.006298"
Set-TextValue "D23" "3.783"
Set-TextValue "D28" "0.0002710"
Set-TextValue "D40" "0.04803"
Set-TextValue "D41" "0.007150"
Set-TextValue "D44" "0.01185"
Set-TextValue "D45" "0.00007080"
Set-TextValue "D47" "0.9346"
Set-TextValue "D48" "0.008390"
Set-TextValue "D49" "0.00001502"
Set-TextValue "D50" "0.01242"

# --- Rows 42/43: BKEXToken and CEJI swap places ----------------------------
# Row 42 was BKEXToken, becomes CEJI; row 43 was CEJI, becomes BKEXToken.
# The rank prefix in column E (41.. / 42..) stays tied to the row, only the
# coin name suffix changes; the Price values are each row's freshly scraped
# number (not simply swapped from the other row).
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004507"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1117"
$ws.Range("E43").Value = "42BKEXTokenBKK"
